$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.783.90"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "3.809.05"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.83"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.11"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.452"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000251"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.04"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "4.447.77"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "3.792.55"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "67.824.58"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.45"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.42"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.89"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000147"
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.32"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.12"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "3.958.36"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("E31").Value = "  +2.61%  "
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.49"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.08"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.994"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.21"
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.28"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.300"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.34"
$ws.Range("E46").Value = "  +8.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.72"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.38"
$ws.Range("E48").Value = "  +11.39%  "
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "390.47"
$ws.Range("E51").Value = "  +0.28%  "
